# Update "想去人数" (interested-count) figures across the workbook's four
# sheets to match the newly scraped totals.
#
# Sheet 1 "展览"     (Exhibitions)
# Sheet 2 "演出"     (Performances)
# Sheet 3 "本地生活" (Local life)
# Sheet 4 "全部类型" (All types - aggregate of the above three)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 9922
$ws1.Range("F13").Value = 1551
$ws1.Range("F18").Value = 454
$ws1.Range("F19").Value = 1139
$ws1.Range("F20").Value = 119
$ws1.Range("F26").Value = 293
$ws1.Range("F31").Value = 23
$ws1.Range("F34").Value = 214
$ws1.Range("F36").Value = 354
$ws1.Range("F38").Value = 458
$ws1.Range("F41").Value = 1261
$ws1.Range("F44").Value = 311

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F6").Value = 81
$ws2.Range("F8").Value = 714
$ws2.Range("F11").Value = 2
$ws2.Range("F20").Value = 441

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 800
$ws3.Range("F5").Value = 177
$ws3.Range("F6").Value = 2473
$ws3.Range("F7").Value = 3939
$ws3.Range("F10").Value = 211

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 800
$ws4.Range("F4").Value = 9922
$ws4.Range("F6").Value = 3939
$ws4.Range("F8").Value = 211
$ws4.Range("F9").Value = 211
$ws4.Range("F11").Value = 1551
$ws4.Range("F16").Value = 454
$ws4.Range("F17").Value = 1139
$ws4.Range("F18").Value = 119
$ws4.Range("F27").Value = 293
$ws4.Range("F32").Value = 23
$ws4.Range("F37").Value = 354
$ws4.Range("F39").Value = 458
$ws4.Range("F47").Value = 311
